$wb = $excel.ActiveWorkbook
$wsOverall = $wb.Worksheets.Item("Overall Test Report")
$wsTests   = $wb.Worksheets.Item("Test Cases & Results")
$wsEnums   = $wb.Worksheets.Item("Enums")

# ---------------------------------------------------------------------------
# "Test Cases & Results" sheet: shift header/table down by one row (was at
# row 2, now at row 3) and fill in the 15 system test cases.
# ---------------------------------------------------------------------------

$wsTests.Rows.Item(2).Insert()

# Data rows (row 4 .. row 18)
$rows = @(
    @{ D="REQ-01"; E="Mid Impact"; F="Test that Web Server pages are shown correctly and content is shown (products)"; G="Webserver Container is running (SPmartWeb)"; H="Open Localhost/IP Adderess of RPI"; I="Products can be seen and added to cart" },
    @{ D="REQ-02"; E="Mid Impact"; F="Test that Web Server page shortcuts are shown correctly "; G="Webserver Container is running (SPmartWeb)"; H="Open Localhost/IP Adderess of RPI"; I="Shortcuts can be seen in Navbar" },
    @{ D="REQ-03"; E="High Impact"; F="Test that cart pages functions properly"; G="Webserver Container is running (SPmartWeb) Products are added to cart and user is logged in"; H="Add a product to cart"; I="Products are shown in cart and checkout options are shown" },
    @{ D="REQ-04"; E="High Impact"; F="Test that checkout options works properly"; G="Webserver Container is running (SPmartWeb) Products are added to cart and user is logged in"; H="Choose a checkout option and checkout"; I="Additional price for delivery is shown and added if selected. When checkout, sent to My Orders" },
    @{ D="REQ-05"; E="Mid Impact"; F="Test that orders are shown and Self Pickup QR codes are shown"; G="Webserver Container is running (SPmartWeb). Order is placed and user is logged in"; H="Checkout/Go to My Orders page"; I="Orders are shown and QR codes shown where applicable" },
    @{ D="REQ-06"; E="Low Impact"; F="Test that LCD shows menu as required"; G="Python Container is running (SPmartIRL)"; H="Run SPmartIRL container"; I="displayed on the LCD Screen`nLine 1: `u{201C}SPmart Menu`u{201D}`nLine 2: `u{201C}1. Self-Checkout`u{201D}" },
    @{ D="REQ-07"; E="High Impact"; F='Test that pressing "1" starts self checkout process'; G="Python Container is running (SPmartIRL)"; H='Press "1" on keypad'; I="REQ-08 menu shown" },
    @{ D="REQ-08"; E="Low Impact"; F="Test that LCD shows instruction as required"; G='Python Container is running (SPmartIRL), "1" Is pressed in main menu'; H="Nil (auto runs)"; I="displayed on the LCD Screen:`nLine 1: `u{201C}Scan item`u{201D}`nLine 2: `u{201C}at camera`u{201D}" },
    @{ D="REQ-09-REQ-17"; E="High Impact"; F="Test that Product Scanning works"; G='Python Container is running (SPmartIRL), "1" Is pressed in main menu'; H="Nil (auto runs)"; I="Flowchart in Figure 1 (SRS Doc) is followed" },
    @{ D="REQ-18"; E="Mid Impact"; F="Test that Payment Menu shown as required"; G='Python Container is running (SPmartIRL), "1" Is pressed in main menu, Products are scanned and user proceeds to payment'; H='Press "1" on keypad after scanning products'; I=" Display the following lines on the LCD `nLine 1: `u{201C}Payment Method?`u{201D}`nLine 2: `u{201C}1. ATM, 2. PayWave`u{201D}"; J="Display the following lines on the LCD `nLine 1: `u{201C}Payment Method?`u{201D}`nLine 2: `u{201C}1. ATM, 2. PayWave`u{201D}" },
    @{ D="REQ-19-REQ-35"; E="High Impact"; F="(PyTest) Test that Payment works"; G='Python Container is running (SPmartIRL), "1" Is pressed in main menu, Products are scanned and user proceeds to payment'; H=$null; I="Flowchart in Figure 2 (SRS Doc) is followed" },
    @{ D="REQ-9"; E="High Impact"; F="(PyTest) Test that Barcode Scanning Works"; G="MySQL Container is running, PyTest Unit Test test_camera_scanning is ran."; H="Point camera at a barcode"; I="LCD to display product details fetched from database" },
    @{ D="REQ-13"; E="High Impact"; F="(PyTest) Test that Python is able to fetch data from MySQL databse"; G="MySQL Container is running, PyTest Unit Test test_db is ran."; H="Run all tests in test_db.py"; I="LCD to display product details fetched from database" },
    @{ D="REQ-23"; E="Mid Impact"; F="(PyTest) Test that entering right/wrong PIN works as intended"; G="Python Container is running, PyTest Unit Test test_verify_pin() is ran."; H="Run test_verify_pin() in test_main.py"; I="Test pass" },
    @{ D="REQ-28-REQ-33"; E="Mid Impact"; F="(PyTest) Test that RFID payment (PayWave) is working as intended"; G="Python Container is running, PyTest Unit Test test_pay_with_paywave  is ran."; H=$null; I="Test pass, Flowchart steps in Figure 2 REQ-28-33 is followed on RPI Dev Board" }
)

$r = 4
foreach ($row in $rows) {
    $wsTests.Range("D$r").Value = $row.D
    $wsTests.Range("E$r").Value = $row.E
    $wsTests.Range("F$r").Value = $row.F
    $wsTests.Range("G$r").Value = $row.G
    if ($null -ne $row.H) {
        $wsTests.Range("H$r").Value = $row.H
    }
    $wsTests.Range("I$r").Value = $row.I
    if ($row.ContainsKey("J")) {
        $wsTests.Range("J$r").Value = $row.J
    } else {
        $wsTests.Range("J$r").Value = $row.I
    }
    $wsTests.Range("K$r").Value = "Pass"
    $r = $r + 1
}

$wsOverall.Range("C3").Formula = "=COUNTIF('Test Cases & Results'!B4:B47, ""<>"")"
$wsOverall.Range("C4").Formula = "=COUNTIF('Test Cases & Results'!K4:K49, ""Pass"")"
$wsOverall.Range("C5").Formula = "=COUNTIF('Test Cases & Results'!K4:K49, ""Fail"")"
$wsOverall.Range("C6").Formula = "=COUNTIF('Test Cases & Results'!K4:K49, ""Not Tested"")"
